# chore: update Sheets via scheduled runner
#
# Refreshes the market-board pricing columns (H:N) on a handful of Leve
# rows across the Carbuncle "Profits" workbook's class sheets. These are
# plain scraped values (no formulas in this workbook) so each cell is set
# directly; a couple of rows gain/lose an M or N cell entirely, matching
# how the source feed only emits a profit cell when the corresponding
# cost column is populated.

$wb = $excel.ActiveWorkbook

# --- ALC ------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 3276.5469
$ws.Range("I76").Value = 3013.9106
$ws.Range("J76").Value = 5115
$ws.Range("K76").Value = 3013.9106
$ws.Range("L76").Value = 5115
$ws.Range("M76").Value = -2698.9106
$ws.Range("N76").Value = -5745

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 3276.5469
$ws.Range("I79").Value = 3013.9106
$ws.Range("J79").Value = 5115
$ws.Range("K79").Value = 3013.9106
$ws.Range("L79").Value = 5115
$ws.Range("M79").Value = -1921.9106
$ws.Range("N79").Value = -7299

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1873.4286
$ws.Range("I137").Value = 1352.2084
$ws.Range("J137").Value = 5000.75
$ws.Range("K137").Value = 4056.6252
$ws.Range("L137").Value = 15002.25
$ws.Range("M137").Value = -1506.6252
$ws.Range("N137").Value = -20102.25

# --- ARM ------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 48: Skillet to the Stars / Steel Frypan
$ws.Range("H48").Value = 80146
$ws.Range("J48").Value = 80146
$ws.Range("L48").Value = 80146
$ws.Range("N48").Value = -80914

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4781.595
$ws.Range("I61").Value = 3811.0605
$ws.Range("J61").Value = 8340.223
$ws.Range("K61").Value = 3811.0605
$ws.Range("L61").Value = 8340.223
$ws.Range("M61").Value = -3599.0605
$ws.Range("N61").Value = -8764.223

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1957.25
$ws.Range("I102").Value = 1943
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1943
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -321
$ws.Range("N102").Value = -5244

# Row 124: Ace of Gloves / High Durium Gauntlets of Fending
$ws.Range("H124").Value = 16095.667
$ws.Range("J124").Value = 16095.667
$ws.Range("L124").Value = 16095.667
$ws.Range("N124").Value = -25915.667

# Row 125: The Incomplete Costume / High Durium Armor of Fending
$ws.Range("H125").Value = 75686.25
$ws.Range("J125").Value = 75686.25
$ws.Range("L125").Value = 75686.25
$ws.Range("N125").Value = -85526.25

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4781.595
$ws.Range("I136").Value = 3811.0605
$ws.Range("J136").Value = 8340.223
$ws.Range("K136").Value = 11433.1815
$ws.Range("L136").Value = 25020.669
$ws.Range("M136").Value = -8883.181500000001
$ws.Range("N136").Value = -30120.669

# --- BSM ------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 41: A Spy in the House of Love / Steel Awl
$ws.Range("H41").Value = 87475
$ws.Range("J41").Value = 87475
$ws.Range("L41").Value = 87475
$ws.Range("N41").Value = -88251

# Row 42: Hard Knock Life / Steel Sledgehammer
$ws.Range("H42").Value = 70342
$ws.Range("J42").Value = 70342
$ws.Range("L42").Value = 70342
$ws.Range("N42").Value = -70998

# Row 47: Lending a Hand / Steel Raising Hammer
$ws.Range("H47").Value = 76861.336
$ws.Range("J47").Value = 76861.336
$ws.Range("L47").Value = 76861.336
$ws.Range("N47").Value = -77901.336

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 25251.691
$ws.Range("I99").Value = 56342.223
$ws.Range("J99").Value = 1933.7916
$ws.Range("K99").Value = 56342.223
$ws.Range("L99").Value = 1933.7916
$ws.Range("M99").Value = -54844.223
$ws.Range("N99").Value = -4929.7916

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 3213.2222
$ws.Range("I105").Value = 3364.875
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 3364.875
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -1617.875
$ws.Range("N105").Value = -5494

# --- CRP ------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 6586.025
$ws.Range("I31").Value = 1205.5834
$ws.Range("J31").Value = 55010
$ws.Range("K31").Value = 1205.5834
$ws.Range("L31").Value = 55010
$ws.Range("M31").Value = -910.5834
$ws.Range("N31").Value = -55600

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 6586.025
$ws.Range("I34").Value = 1205.5834
$ws.Range("J34").Value = 55010
$ws.Range("K34").Value = 1205.5834
$ws.Range("L34").Value = 55010
$ws.Range("M34").Value = -1003.5834
$ws.Range("N34").Value = -55414

# Row 62: Splinter in the Sewers / Cedar Lumber (M62 drops out)
$ws.Range("H62").Value = 3300.7144
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3300.7144
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3300.7144
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4548.7144

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber (M65 drops out)
$ws.Range("H65").Value = 3300.7144
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3300.7144
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 16503.572
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -22743.572

# --- CUL ------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1170.5
$ws.Range("I5").Value = 616.8333
$ws.Range("J5").Value = 1645.0714
$ws.Range("K5").Value = 1850.4999
$ws.Range("L5").Value = 4935.2142
$ws.Range("M5").Value = -1738.4999
$ws.Range("N5").Value = -5159.2142

# Row 61: Red Letter Day / Rolanberry Lassi
$ws.Range("H61").Value = 307.3
$ws.Range("I61").Value = 245.6
$ws.Range("J61").Value = 369
$ws.Range("K61").Value = 736.8
$ws.Range("L61").Value = 1107
$ws.Range("M61").Value = -521.8
$ws.Range("N61").Value = -1537

# Row 122: Salt of the North / Northern Sea Salt (M122 newly added)
$ws.Range("H122").Value = 1329.6666
$ws.Range("I122").Value = 995
$ws.Range("K122").Value = 8955
$ws.Range("M122").Value = -6505

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1022.75
$ws.Range("I132").Value = 878.6667
$ws.Range("J132").Value = 1455
$ws.Range("K132").Value = 7908.0003
$ws.Range("L132").Value = 13095
$ws.Range("M132").Value = -5378.0003
$ws.Range("N132").Value = -18155

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1170.5
$ws.Range("I135").Value = 616.8333
$ws.Range("J135").Value = 1645.0714
$ws.Range("K135").Value = 5551.4997
$ws.Range("L135").Value = 14805.6426
$ws.Range("M135").Value = -3016.4997
$ws.Range("N135").Value = -19875.6426

# --- GSM ------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 5260.25
$ws.Range("I97").Value = 3676.6667
$ws.Range("J97").Value = 10011
$ws.Range("K97").Value = 3676.6667
$ws.Range("L97").Value = 10011
$ws.Range("M97").Value = -3180.6667
$ws.Range("N97").Value = -11003

# --- LTW ------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 100: Tiger in the Sack / Tiger Leather (M100 & N100 drop out)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3311.2666
$ws.Range("I122").Value = 2168.3
$ws.Range("J122").Value = 5597.2
$ws.Range("K122").Value = 6504.900000000001
$ws.Range("L122").Value = 16791.6
$ws.Range("M122").Value = -4054.900000000001
$ws.Range("N122").Value = -21691.6

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 6588.625
$ws.Range("I136").Value = 4784
$ws.Range("K136").Value = 14352
$ws.Range("M136").Value = -11802

# --- WVR ------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 10453208
$ws.Range("I136").Value = 47762620
$ws.Range("J136").Value = 6572.8
$ws.Range("K136").Value = 143287860
$ws.Range("L136").Value = 19718.4
$ws.Range("M136").Value = -143285310
$ws.Range("N136").Value = -24818.4
